# Key Features Og AngularJS
# Recolor the "7 | Custom Directives | Properties of Custom Directives |
# 4 hrs | 2 days" syllabus row (table row 11, 1-based) from the default
# black/theme-text1 color to green (RGB 0x00B050), matching the author's
# highlight of that row in the source diff.

$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

$greenColor = 5287936  # RGB(0, 176, 80) == hex 00B050

$targetFirstCell = "Custom Directives"
$targetRow = $null

for ($r = 1; $r -le $table.Rows.Count; $r++) {
    $row = $table.Rows.Item($r)
    $secondCellText = $row.Cells.Item(2).Range.Text.Trim()
    if ($secondCellText -eq $targetFirstCell) {
        $targetRow = $row
        break
    }
}

if ($targetRow -eq $null) {
    # Fallback: row 11 is the "7 / Custom Directives / ..." row in the
    # known document layout.
    $targetRow = $table.Rows.Item(11)
}

for ($c = 1; $c -le $targetRow.Cells.Count; $c++) {
    $cell = $targetRow.Cells.Item($c)
    $cell.Range.Font.Color = $greenColor
}

Write-Host "Recolored row to RGB 00B050"
